# Commit final, todas as cases modificação com os padrões Factory
#
# Rebuilds "Pesquisa" (sheet2) as a small product catalog with a new
# "Descricao" column (Sucesso/Falha rows), and updates "Cadastro" (sheet1):
# the logged-in userName, the Pais/Brazil pair, and column A's width.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Sheet2 ("Pesquisa") — Phase A: capture/apply formatting BEFORE values are
# cleared (copy sources must be read before they get overwritten).
# ---------------------------------------------------------------------------

# B4 (old "empty note" cell, style carries a plain font) -> new F9 marker cell.
$ws2.Range("B4").Copy()
$ws2.Range("F9").PasteSpecial($xlPasteFormats)

# A2 (old centered/filled style) -> the new product-category column A2:A5.
$ws2.Range("A2").Copy()
$ws2.Range("A2:A5").PasteSpecial($xlPasteFormats)

# I15 (plain data style) -> the new description columns B/C (this also
# overwrites B4's old style, which is fine - already captured above).
$ws2.Range("I15").Copy()
$ws2.Range("B2:C5").PasteSpecial($xlPasteFormats)
$ws2.Range("B7:C8").PasteSpecial($xlPasteFormats)

# A1 (header style) -> the full new header row A1:C1.
$ws2.Range("A1").Copy()
$ws2.Range("A1:C1").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# H12 already carries the right style and position - left untouched.
# I15 was only needed as a formatting donor; it does not survive into the new
# layout, so drop it (value + style) entirely.
$ws2.Range("I15").Clear()

# Phase B: wipe the old cell values only, keeping the formatting just laid down.
$ws2.Cells.ClearContents()

# ---------------------------------------------------------------------------
# Sheet2 — Phase C: write the new values. Order matters: it must match the
# sequence new strings are first introduced so the shared-string table comes
# out the same as a natural top-to-bottom data-entry pass would produce.
# ---------------------------------------------------------------------------
$ws2.Range("B2").Value = "Bose Soundlink Bluetooth Speaker III"
$ws2.Range("A4").Value = "Mice"
$ws2.Range("C1").Value = "Descricao"
$ws2.Range("C2").Value = "Sucesso"
$ws2.Range("C7").Value = "Falha"
$ws2.Range("B3").Value = "Bose Soundlink Bluetooth Speaker "
$ws2.Range("A5").Value = "Laptops"
$ws2.Range("B5").Value = "HP Chromebook 14 G1(ENERGY STAR)"
$ws2.Range("B7").Value = "Banana"
$ws2.Range("B8").Value = "Mochila"

# Remaining cells reuse already-existing shared strings.
$ws2.Range("A1").Value = "Produto"
$ws2.Range("B1").Value = "nomeProduto"
$ws2.Range("A2").Value = "Speakers"
$ws2.Range("A3").Value = "Speakers"
$ws2.Range("C3").Value = "Sucesso"
$ws2.Range("B4").Value = "HP USB 3 Button Optical Mouse"
$ws2.Range("C4").Value = "Sucesso"
$ws2.Range("C5").Value = "Sucesso"
$ws2.Range("C8").Value = "Falha"

# ---------------------------------------------------------------------------
# Sheet1 ("Cadastro") — update userName + Pais/Brazil, widen column A.
# ---------------------------------------------------------------------------
$ws1.Range("K1").Value = "Pais"
$ws1.Range("K2").Value = "Brazil"
$ws1.Range("A2").Value = "lucasViado"

$ws1.Columns.Item(1).ColumnWidth = 13.6666667

# ---------------------------------------------------------------------------
# Restore the on-screen selection: Pesquisa's cursor ends on M11, Cadastro
# stays the active tab (its own A2 selection is untouched).
# ---------------------------------------------------------------------------
[void]$ws2.Activate()
[void]$ws2.Range("M11").Select()
[void]$ws1.Activate()
